# Updates cell values on Sheet1 to match the scraped coinranking.com
# snapshot for this commit. Values that look numeric (Price / Volume(1h))
# are forced to text with a leading apostrophe so they keep their exact
# printed form (trailing zeros, %, thousands separators) exactly like the
# original inline-string cells, then Style is reset to Normal so no stray
# quote-prefix formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''245.38'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '''-0.17%'
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = '''26.51'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '''4.22%'
$ws.Range('E3').Style = "Normal"
$ws.Range('D4').Value = '''5.144'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '''0.43%'
$ws.Range('E4').Style = "Normal"
$ws.Range('E5').Value = '''-0.07%'
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = '''6.475'
$ws.Range('D6').Style = "Normal"
$ws.Range('E7').Value = '''0.06%'
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = '''0.8397'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '''-0.19%'
$ws.Range('E8').Style = "Normal"
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').Value = '''0.1330'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '''-0.97%'
$ws.Range('E9').Style = "Normal"
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D10').Value = '''0.06993'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '''0.68%'
$ws.Range('E10').Style = "Normal"
$ws.Range('B11').Value = 'BitrueCoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D11').Value = '''0.02869'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '''0.40%'
$ws.Range('E11').Style = "Normal"
$ws.Range('B12').Value = 'BitMartToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D12').Value = '''0.09381'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '''-0.06%'
$ws.Range('E12').Style = "Normal"
$ws.Range('B13').Value = 'BitForexToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D13').Value = '''0.001528'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '''0.20%'
$ws.Range('E13').Style = "Normal"
$ws.Range('B14').Value = 'One'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D14').Value = '''0.0005982'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '''0.51%'
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = '''0.006131'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '''-0.05%'
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = '''3.634'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '''3.48%'
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = '''3.037'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '''0.68%'
$ws.Range('E17').Style = "Normal"
$ws.Range('E18').Value = '''4.85%'
$ws.Range('E18').Style = "Normal"
$ws.Range('E19').Value = '''-2.11%'
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = '''0.03066'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '''-3.16%'
$ws.Range('E20').Style = "Normal"
$ws.Range('E21').Value = '''-2.20%'
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = '''3.750'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '''-0.02%'
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = '''0.04596'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '''-2.61%'
$ws.Range('E23').Style = "Normal"
$ws.Range('E24').Value = '''2.47%'
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = '''0.001248'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '''-0.07%'
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = '''0.004514'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '''5.83%'
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = '''0.00009599'
$ws.Range('D27').Style = "Normal"
$ws.Range('D28').Value = '''0.0001397'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '''0.59%'
$ws.Range('E28').Style = "Normal"
$ws.Range('D40').Value = '''0.03641'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '''-0.64%'
$ws.Range('E40').Style = "Normal"
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D41').Value = '''0.1376'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '''1.38%'
$ws.Range('E41').Style = "Normal"
$ws.Range('B42').Value = 'CEJI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D42').Value = '''0.002550'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '''-2.86%'
$ws.Range('E42').Style = "Normal"
$ws.Range('B43').Value = 'KickToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D43').Value = '''0.003452'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '''-44.70%'
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = '''0.008212'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '''-2.62%'
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = '''0.00005359'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '''1.08%'
$ws.Range('E45').Style = "Normal"
$ws.Range('E46').Value = '''0.01%'
$ws.Range('E46').Style = "Normal"
$ws.Range('E47').Value = '''-51.55%'
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = '''0.002550'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '''20.47%'
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = '''0.00002101'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '''0.01%'
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = '''0.0002001'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '''0.01%'
$ws.Range('E50').Style = "Normal"
